# Auto-generated Excel COM-interop script
# Applies scheduled-runner profit/price updates to the Cactuar_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 552.3077
$ws.Range("I5").Value = 182.08333
$ws.Range("K5").Value = 182.08333
$ws.Range("M5").Value = -67.08332999999999
$ws.Range("H51").Value = 9248.75
$ws.Range("J51").Value = 9248.75
$ws.Range("L51").Value = 9248.75
$ws.Range("N51").Value = -10216.75
$ws.Range("H74").Value = 6999.5
$ws.Range("I74").Value = 6666
$ws.Range("K74").Value = 6666
$ws.Range("M74").Value = -5730
$ws.Range("H77").Value = 6999.5
$ws.Range("I77").Value = 6666
$ws.Range("K77").Value = 33330
$ws.Range("M77").Value = -28650
$ws.Range("H121").Value = 3986.8076
$ws.Range("J121").Value = 3986.8076
$ws.Range("L121").Value = 11960.4228
$ws.Range("N121").Value = -15454.4228
$ws.Range("H135").Value = 2734.1177
$ws.Range("I135").Value = 1617.75
$ws.Range("J135").Value = 5413.4
$ws.Range("K135").Value = 14559.75
$ws.Range("L135").Value = 48720.6
$ws.Range("M135").Value = -12024.75
$ws.Range("N135").Value = -53790.6
$ws.Range("H138").Value = 5480.5586
$ws.Range("I138").Value = 1598.2222
$ws.Range("J138").Value = 6878.2
$ws.Range("K138").Value = 4794.6666
$ws.Range("L138").Value = 20634.6
$ws.Range("M138").Value = 345.3334000000004
$ws.Range("N138").Value = -30914.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5570.381
$ws.Range("I61").Value = 5471.919
$ws.Range("K61").Value = 5471.919
$ws.Range("M61").Value = -5259.919
$ws.Range("H74").Value = 8335035.5
$ws.Range("I74").Value = 8930074
$ws.Range("K74").Value = 8930074
$ws.Range("M74").Value = -8929200
$ws.Range("H77").Value = 8335035.5
$ws.Range("I77").Value = 8930074
$ws.Range("K77").Value = 44650370
$ws.Range("M77").Value = -44646002
$ws.Range("H101").Value = 125000
$ws.Range("J101").Value = 125000
$ws.Range("L101").Value = 125000
$ws.Range("N101").Value = -131490
$ws.Range("H135").Value = 99997.664
$ws.Range("J135").Value = 99997.664
$ws.Range("L135").Value = 99997.664
$ws.Range("N135").Value = -110137.664
$ws.Range("H136").Value = 5570.381
$ws.Range("I136").Value = 5471.919
$ws.Range("K136").Value = 16415.757
$ws.Range("M136").Value = -13865.757

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1035.4286
$ws.Range("I80").Value = 396
$ws.Range("K80").Value = 396
$ws.Range("M80").Value = 602
$ws.Range("H83").Value = 1035.4286
$ws.Range("I83").Value = 396
$ws.Range("K83").Value = 1980
$ws.Range("M83").Value = 3012
$ws.Range("H94").Value = 2163.1365
$ws.Range("I94").Value = 1896.2667
$ws.Range("J94").Value = 2735
$ws.Range("K94").Value = 1896.2667
$ws.Range("L94").Value = 2735
$ws.Range("M94").Value = -1445.2667
$ws.Range("N94").Value = -3637
$ws.Range("H124").Value = 29285.715
$ws.Range("J124").Value = 29285.715
$ws.Range("L124").Value = 29285.715
$ws.Range("N124").Value = -39105.715
$ws.Range("H134").Value = 3255.4
$ws.Range("I134").Value = 3132.7083
$ws.Range("K134").Value = 9398.124899999999
$ws.Range("M134").Value = -6863.124899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 749.3
$ws.Range("J22").Value = 824.75
$ws.Range("L22").Value = 824.75
$ws.Range("N22").Value = -1524.75
$ws.Range("H105").Value = 1413.5714
$ws.Range("I105").Value = 1223.75
$ws.Range("K105").Value = 1223.75
$ws.Range("M105").Value = 523.25
$ws.Range("H132").Value = 47621536
$ws.Range("I132").Value = 58826084
$ws.Range("J132").Value = 2209.75
$ws.Range("K132").Value = 176478252
$ws.Range("L132").Value = 6629.25
$ws.Range("M132").Value = -176475722
$ws.Range("N132").Value = -11689.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 422
$ws.Range("J33").Value = 131.25
$ws.Range("L33").Value = 787.5
$ws.Range("N33").Value = -1353.5
$ws.Range("H61").Value = 374.1
$ws.Range("I61").Value = 105.166664
$ws.Range("K61").Value = 315.499992
$ws.Range("M61").Value = -100.499992
$ws.Range("H131").Value = 16826328
$ws.Range("J131").Value = 20478910
$ws.Range("L131").Value = 61436730
$ws.Range("N131").Value = -61446810

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2531.1
$ws.Range("I97").Value = 2256.889
$ws.Range("J97").Value = 4999
$ws.Range("K97").Value = 2256.889
$ws.Range("L97").Value = 4999
$ws.Range("M97").Value = -1760.889
$ws.Range("N97").Value = -5991
$ws.Range("H102").Value = 12446298
$ws.Range("I102").Value = 17593672
$ws.Range("K102").Value = 17593672
$ws.Range("M102").Value = -17592050
$ws.Range("H113").Value = 1450
$ws.Range("I113").Value = 1400
$ws.Range("K113").Value = 1400
$ws.Range("M113").Value = 770
$ws.Range("H132").Value = 47406.74
$ws.Range("I132").Value = 60897.2
$ws.Range("K132").Value = 182691.6
$ws.Range("M132").Value = -180161.6
$ws.Range("H138").Value = 77777
$ws.Range("J138").Value = 77777
$ws.Range("L138").Value = 77777
$ws.Range("N138").Value = -88057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6571.4287
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888
$ws.Range("H126").Value = 6571.4287
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 5146.83
$ws.Range("I132").Value = 4593.5693
$ws.Range("K132").Value = 13780.7079
$ws.Range("M132").Value = -11250.7079
$ws.Range("H136").Value = 6914.2856
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 7400
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 22200
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = -27300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 47999.4
$ws.Range("J94").Value = 47999.4
$ws.Range("L94").Value = 47999.4
$ws.Range("N94").Value = -49801.4
$ws.Range("H132").Value = 4763
$ws.Range("I132").Value = 8997
$ws.Range("J132").Value = 4570.5454
$ws.Range("K132").Value = 26991
$ws.Range("L132").Value = 13711.6362
$ws.Range("M132").Value = -24461
$ws.Range("N132").Value = -18771.6362
$ws.Range("H136").Value = 5075.7744
$ws.Range("I136").Value = 2729.6924
$ws.Range("J136").Value = 9053.913
$ws.Range("K136").Value = 8189.0772
$ws.Range("L136").Value = 27161.739
$ws.Range("M136").Value = -5639.0772
$ws.Range("N136").Value = -32261.739

Write-Host "Applied 168 cell updates across 8 sheets."
